{"js": "// Generalize \"CivicActions\" department references in the PL section of the\n// System Security Plan, and fix a few grammar nits, per the commit:\n//   \"New 'Contractor' component (generalized 'CivicActions'); use secrender\n//    against templates\"\n//\n// Each replacement below targets the full, unique text of one paragraph run\n// (verified against the source OOXML) so a plain text search/replace is\n// unambiguous.\n\nconst replacements = [\n  [\n    \"The SSP is reviewed and approved by the authorizing official prior to plan implementation. A copy of the SSP is provided to authorized CivicActions and assessing personnel including the System Owner, Authorizing Official, Information System Security Officer, System/Network Administrator and CivicActions Operations. The SSP is maintained by CivicActions Security.\",\n    \"The SSP is reviewed and approved by the authorizing official prior to plan implementation. A copy of the SSP is provided to authorized CivicActions and assessing personnel including the System Owner, Authorizing Official, Information System Security Officer, System/Network Administrator, and CivicActions\\u2019 Operations staff. The SSP is maintained by CivicActions\\u2019 Security Office.\"\n  ],\n  [\n    \"The SSP is reviewed at least annually by the System Owner and CivicActions Operations in collaboration with CivicActions Security.\",\n    \"The SSP is reviewed at least annually by the System Owner and CivicActions\\u2019 Operations staff in collaboration with CivicActions\\u2019 Security Office.\"\n  ],\n  [\n    \"CivicActions Operations in collaboration with CivicActions Security updates the system description and control descriptions within the SSP as needed to verify the SSP is an accurate description of the system.\",\n    \"CivicActions\\u2019 Operations staff in collaboration with CivicActions\\u2019 Security Office updates the system description and control descriptions within the SSP as needed to verify the SSP is an accurate description of the system.\"\n  ],\n  [\n    \"CivicActions has created and made readily available to individuals requiring access to the information system the rules that describes their responsibilities and expected behavior with regard to information and information system usage. These rules, defined as the Acceptable Use Policy, are included in the CivicActions Security Policy accessible here :\",\n    \"CivicActions has created and made readily available to individuals requiring access to the information system the rules that describe their responsibilities and expected behavior with regard to information and information system usage. These rules, defined as the Acceptable Use Policy, are included in the CivicActions Security Policy accessible here:\"\n  ],\n  [\n    \"CivicActions HR receives a signed acknowledgment from all employees, indicating that they have read, understand, and agree to abide by the rules of behavior, before authorizing access to information and the information system. The text of the electronically signed (via DocuSign) acknowledgement document has been uploaded to CSAM as artifact:\",\n    \"CivicActions HR receives a signed acknowledgment from all employees, indicating that they have read, understand, and agree to abide by the rules of behavior, before authorizing access to information and the information system. The text of the electronically signed (via DocuSign) acknowledgment document has been uploaded to CSAM as artifact:\"\n  ],\n  [\n    \"CivicActions reviews the CivicActions Security Policy at least annually and updates is as required.\",\n    \"CivicActions reviews the CivicActions Security Policy at least annually and updates as required.\"\n  ],\n  [\n    \"CivicActions requires individuals who have signed a previous version of the CivicActions Security Policy to read and re-sign when any part of it, including the Acceptable Use Policy/Rules of Behavior, are revised/updated.\",\n    \"CivicActions requires individuals who have signed a previous version of the CivicActions Security Policy to read and re-sign when any part of it, including the Acceptable Use Policy/Rules of Behavior, is revised/updated.\"\n  ]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Generalize \"CivicActions\" department references in the PL section of the\n# System Security Plan, and fix a few grammar nits, per the commit:\n#   \"New 'Contractor' component (generalized 'CivicActions'); use secrender\n#    against templates\"\n#\n# Each replacement below targets the full, unique text of one paragraph run\n# (verified against the source OOXML) so a plain Find/Replace is unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\n    \"The SSP is reviewed and approved by the authorizing official prior to plan implementation. A copy of the SSP is provided to authorized CivicActions and assessing personnel including the System Owner, Authorizing Official, Information System Security Officer, System/Network Administrator and CivicActions Operations. The SSP is maintained by CivicActions Security.\",\n    \"The SSP is reviewed and approved by the authorizing official prior to plan implementation. A copy of the SSP is provided to authorized CivicActions and assessing personnel including the System Owner, Authorizing Official, Information System Security Officer, System/Network Administrator, and CivicActions\u2019 Operations staff. The SSP is maintained by CivicActions\u2019 Security Office.\"\n  ),\n  @(\n    \"The SSP is reviewed at least annually by the System Owner and CivicActions Operations in collaboration with CivicActions Security.\",\n    \"The SSP is reviewed at least annually by the System Owner and CivicActions\u2019 Operations staff in collaboration with CivicActions\u2019 Security Office.\"\n  ),\n  @(\n    \"CivicActions Operations in collaboration with CivicActions Security updates the system description and control descriptions within the SSP as needed to verify the SSP is an accurate description of the system.\",\n    \"CivicActions\u2019 Operations staff in collaboration with CivicActions\u2019 Security Office updates the system description and control descriptions within the SSP as needed to verify the SSP is an accurate description of the system.\"\n  ),\n  @(\n    \"CivicActions has created and made readily available to individuals requiring access to the information system the rules that describes their responsibilities and expected behavior with regard to information and information system usage. These rules, defined as the Acceptable Use Policy, are included in the CivicActions Security Policy accessible here :\",\n    \"CivicActions has created and made readily available to individuals requiring access to the information system the rules that describe their responsibilities and expected behavior with regard to information and information system usage. These rules, defined as the Acceptable Use Policy, are included in the CivicActions Security Policy accessible here:\"\n  ),\n  @(\n    \"CivicActions HR receives a signed acknowledgment from all employees, indicating that they have read, understand, and agree to abide by the rules of behavior, before authorizing access to information and the information system. The text of the electronically signed (via DocuSign) acknowledgement document has been uploaded to CSAM as artifact:\",\n    \"CivicActions HR receives a signed acknowledgment from all employees, indicating that they have read, understand, and agree to abide by the rules of behavior, before authorizing access to information and the information system. The text of the electronically signed (via DocuSign) acknowledgment document has been uploaded to CSAM as artifact:\"\n  ),\n  @(\n    \"CivicActions reviews the CivicActions Security Policy at least annually and updates is as required.\",\n    \"CivicActions reviews the CivicActions Security Policy at least annually and updates as required.\"\n  ),\n  @(\n    \"CivicActions requires individuals who have signed a previous version of the CivicActions Security Policy to read and re-sign when any part of it, including the Acceptable Use Policy/Rules of Behavior, are revised/updated.\",\n    \"CivicActions requires individuals who have signed a previous version of the CivicActions Security Policy to read and re-sign when any part of it, including the Acceptable Use Policy/Rules of Behavior, is revised/updated.\"\n  )\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $find = $rng.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 0\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute(\n    $oldText,      # FindText\n    $true,         # MatchCase\n    $false,        # MatchWholeWord\n    $false,        # MatchWildcards\n    $false,        # MatchSoundsLike\n    $false,        # MatchAllWordForms\n    $true,         # Forward\n    0,             # Wrap (wdFindStop)\n    $false,        # Format\n    $newText,      # ReplaceWith\n    2              # Replace (wdReplaceAll)\n  )\n}\n"}
